$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 239, shifting existing rows 239-243 down to 240-244.
$ws.Rows.Item(239).Insert()

# Populate the newly inserted row 239 with the new weekly data point.
$ws.Range("A239").Value = 4
$ws.Range("B239").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C239").Value = "Los Lagos"
$ws.Range("D239").Value = 44448
$ws.Range("E239").Value = 10
$ws.Range("F239").Value = 100112006
$ws.Range("G239").Value = "Repollo"
$ws.Range("H239").Value = "Crespo record"
$ws.Range("I239").Value = "Primera"
$ws.Range("J239").Value = 500
$ws.Range("K239").Value = 1200
$ws.Range("L239").Value = 1200
$ws.Range("M239").Value = 1200
$ws.Range("N239").Value = "$/unidad"
$ws.Range("O239").Value = "Región Metropolitana"
$ws.Range("P239").Value = 1200
$ws.Range("Q239").Value = 1
$ws.Range("R239").Value = "Hortaliza"
